$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 data (MCH201-1) ---
$ws.Range("A2").Value = "MCH201-1"
$ws.Range("C2").Value = "PAPERS, NEWSCUTTINGS, STICKERS POSTERS"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: 24B | GRAP COUNT NUMER: NONE"

# --- Row 3 data (MCH201-2) ---
$ws.Range("A3").Value = "MCH201-2"
$ws.Range("C3").Value = "SLIDES, ARTEFACTS"
$ws.Range("E3").Value = "Series"
$ws.Range("F3").Value = "1 Box"
$ws.Range("G3").Value = "LOCATION: 24B | GRAP COUNT NUMER: NONE"

# --- Empty but formatted cells (D2, H2, D3, H3) ---
# touched implicitly by the font formatting below

# --- Formatting: body font is Calibri 10pt, theme text color ---
$r1 = $ws.Range("A2:A3")
$r1.Font.Name = "Calibri"
$r1.Font.Size = 10
$r1.Font.ThemeColor = 1

$r2 = $ws.Range("C2:E3")
$r2.Font.Name = "Calibri"
$r2.Font.Size = 10
$r2.Font.ThemeColor = 1

$r3 = $ws.Range("G2:H3")
$r3.Font.Name = "Calibri"
$r3.Font.Size = 10
$r3.Font.ThemeColor = 1

$r4 = $ws.Range("F2:F3")
$r4.Font.Name = "Calibri"
$r4.Font.Size = 10
$r4.Font.ThemeColor = 1
$r4.HorizontalAlignment = -4108

# --- Row heights to match the rest of the sheet ---
$ws.Range("A2:A3").RowHeight = 15.75

# --- Restore frozen header pane / selection over the new data range ---
$ws.Range("A2").Select()
$win = $excel.ActiveWindow
$win.FreezePanes = $true
$ws.Range("A2:H3").Select()
